$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.434937333333333
$ws.Range("H2").Value = 4.304812
$ws.Range("I2").Value = 0.5010808920723563
$ws.Range("J2").Value = 0.5010808920723562
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.042868
$ws.Range("N2").Value = 0.128604
$ws.Range("O2").Value = 0.03014606792405771
$ws.Range("P2").Value = 0.03014606792405771
$ws.Range("Q2").Value = 0.06151289360533332
$ws.Range("R2").Value = 0.553616042448
$ws.Range("S2").Value = 0.01510561860786068
$ws.Range("T2").Value = 0.01510561860786068
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.434937333333333
$ws.Range("H3").Value = 4.304812
$ws.Range("I3").Value = 0.5010808920723563
$ws.Range("J3").Value = 0.5010808920723562
$ws.Range("O3").Value = 0.2718481285523376
$ws.Range("P3").Value = 0.2718481285523376
$ws.Range("Q3").Value = 0.5547046815715555
$ws.Range("R3").Value = 4.992342134144001
$ws.Range("S3").Value = 0.1362179027632059
$ws.Range("T3").Value = 0.1362179027632059
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.434937333333333
$ws.Range("H4").Value = 4.304812
$ws.Range("I4").Value = 0.5010808920723563
$ws.Range("J4").Value = 0.5010808920723562
$ws.Range("M4").Value = 0.9839956666666666
$ws.Range("N4").Value = 2.951987
$ws.Range("O4").Value = 0.69197537100662
$ws.Range("P4").Value = 0.69197537100662
$ws.Range("Q4").Value = 1.411972117938222
$ws.Range("R4").Value = 12.707749061444
$ws.Range("S4").Value = 0.3467356361960969
$ws.Range("T4").Value = 0.3467356361960968
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.434937333333333
$ws.Range("H5").Value = 4.304812
$ws.Range("I5").Value = 0.5010808920723563
$ws.Range("J5").Value = 0.5010808920723562
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.008575333333333332
$ws.Range("N5").Value = 0.025726
$ws.Range("O5").Value = 0.006030432516984765
$ws.Range("P5").Value = 0.006030432516984765
$ws.Range("Q5").Value = 0.01230506594577778
$ws.Range("R5").Value = 0.110745593512
$ws.Range("S5").Value = 0.00302173450519287
$ws.Range("T5").Value = 0.00302173450519287
$ws.Range("G6").Value = 0.9964423333333334
$ws.Range("I6").Value = 0.3479582011609289
$ws.Range("J6").Value = 0.3479582011609288
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.042868
$ws.Range("N6").Value = 0.128604
$ws.Range("O6").Value = 0.03014606792405771
$ws.Range("P6").Value = 0.03014606792405771
$ws.Range("Q6").Value = 0.04271548994533333
$ws.Range("R6").Value = 0.384439409508
$ws.Range("S6").Value = 0.0104895715669303
$ws.Range("T6").Value = 0.0104895715669303
$ws.Range("G7").Value = 0.9964423333333334
$ws.Range("I7").Value = 0.3479582011609289
$ws.Range("J7").Value = 0.3479582011609288
$ws.Range("O7").Value = 0.2718481285523376
$ws.Range("P7").Value = 0.2718481285523376
$ws.Range("Q7").Value = 0.3851953770915556
$ws.Range("R7").Value = 3.466758393824001
$ws.Range("S7").Value = 0.09459178580003634
$ws.Range("T7").Value = 0.09459178580003631
$ws.Range("G8").Value = 0.9964423333333334
$ws.Range("I8").Value = 0.3479582011609289
$ws.Range("J8").Value = 0.3479582011609288
$ws.Range("M8").Value = 0.9839956666666666
$ws.Range("N8").Value = 2.951987
$ws.Range("O8").Value = 0.69197537100662
$ws.Range("P8").Value = 0.69197537100662
$ws.Range("Q8").Value = 0.9804949380832223
$ws.Range("R8").Value = 8.824454442749001
$ws.Range("S8").Value = 0.2407785053431299
$ws.Range("T8").Value = 0.2407785053431298
$ws.Range("G9").Value = 0.9964423333333334
$ws.Range("I9").Value = 0.3479582011609289
$ws.Range("J9").Value = 0.3479582011609288
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.008575333333333332
$ws.Range("N9").Value = 0.025726
$ws.Range("O9").Value = 0.006030432516984765
$ws.Range("P9").Value = 0.006030432516984765
$ws.Range("Q9").Value = 0.008544825155777779
$ws.Range("R9").Value = 0.076903426402
$ws.Range("S9").Value = 0.002098338450832391
$ws.Range("T9").Value = 0.002098338450832391
$ws.Range("G10").Value = 0.4323043333333333
$ws.Range("H10").Value = 1.296913
$ws.Range("I10").Value = 0.150960906766715
$ws.Range("J10").Value = 0.1509609067667149
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.042868
$ws.Range("N10").Value = 0.128604
$ws.Range("O10").Value = 0.03014606792405771
$ws.Range("P10").Value = 0.03014606792405771
$ws.Range("Q10").Value = 0.01853202216133333
$ws.Range("R10").Value = 0.166788199452
$ws.Range("S10").Value = 0.004550877749266733
$ws.Range("T10").Value = 0.004550877749266732
$ws.Range("G11").Value = 0.4323043333333333
$ws.Range("H11").Value = 1.296913
$ws.Range("I11").Value = 0.150960906766715
$ws.Range("J11").Value = 0.1509609067667149
$ws.Range("O11").Value = 0.2718481285523376
$ws.Range("P11").Value = 0.2718481285523376
$ws.Range("Q11").Value = 0.1671161743395556
$ws.Range("R11").Value = 1.504045569056
$ws.Range("S11").Value = 0.04103843998909538
$ws.Range("T11").Value = 0.04103843998909536
$ws.Range("G12").Value = 0.4323043333333333
$ws.Range("H12").Value = 1.296913
$ws.Range("I12").Value = 0.150960906766715
$ws.Range("J12").Value = 0.1509609067667149
$ws.Range("M12").Value = 0.9839956666666666
$ws.Range("N12").Value = 2.951987
$ws.Range("O12").Value = 0.69197537100662
$ws.Range("P12").Value = 0.69197537100662
$ws.Range("Q12").Value = 0.4253855906812222
$ws.Range("R12").Value = 3.828470316131
$ws.Range("S12").Value = 0.1044612294673934
$ws.Range("T12").Value = 0.1044612294673933
$ws.Range("G13").Value = 0.4323043333333333
$ws.Range("H13").Value = 1.296913
$ws.Range("I13").Value = 0.150960906766715
$ws.Range("J13").Value = 0.1509609067667149
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.008575333333333332
$ws.Range("N13").Value = 0.025726
$ws.Range("O13").Value = 0.006030432516984765
$ws.Range("P13").Value = 0.006030432516984765
$ws.Range("Q13").Value = 0.003707153759777778
$ws.Range("R13").Value = 0.033364383838
$ws.Range("S13").Value = 0.0009103595609595034
$ws.Range("T13").Value = 0.0009103595609595031

Write-Output "Applied 154 cell updates"
